$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J34").Value = 'AC Power kVA Standard'
$ws.Range("J35").Value = 'Metered Power Charges kVA Standard'
$ws.Range("J36").Value = 'Power Cord kVA Standard'
$ws.Range("D37").Value = 12
$ws.Range("E37").Value = 12
$ws.Range("J37").Value = 'AC Power Monthly Charge'
$ws.Range("D38").Value = 150
$ws.Range("E38").Value = 150
$ws.Range("J38").Value = 'Smart Hands NRC'
$ws.Range("J39").Value = 'Cabinet Installation'
$ws.Range("J40").Value = 'Cage Installation'
$ws.Range("J41").Value = 'Cross Connect Single-Mode Fiber'
$ws.Range("D42").Value = 200
$ws.Range("E42").Value = 200
$ws.Range("J42").Value = 'Equinix Precision Time Standard NTP'
$ws.Range("D43").Value = 200
$ws.Range("E43").Value = 200
$ws.Range("J43").Value = 'Secure Cabinet kVA Standard'
$ws.Range("D44").Value = 200
$ws.Range("E44").Value = 200
$ws.Range("J44").Value = 'AC Circuit 30A 208V'
$ws.Range("D45").Value = 200
$ws.Range("E45").Value = 200
$ws.Range("J45").Value = 'Draw Cap Overage Charge kVA Standard'
$ws.Range("J46").Value = 'Extended Cross Connect Multi-Mode Fiber'
$ws.Range("J47").Value = 'DC Circuit 30 208'
$ws.Range("D48").Value = 12
$ws.Range("E48").Value = 12
$ws.Range("J48").Value = 'AC Power kVA Standard'
$ws.Range("D49").Value = 12
$ws.Range("E49").Value = 12
$ws.Range("J49").Value = 'Metered Power Charges kVA Standard'
$ws.Range("J50").Value = 'Power Cord kVA Standard'
$ws.Range("J51").Value = 'AC Power Monthly Charge'
$ws.Range("D52").Value = 150
$ws.Range("E52").Value = 150
$ws.Range("J52").Value = 'Smart Hands NRC'
$ws.Range("D53").Value = 200
$ws.Range("E53").Value = 200
$ws.Range("J53").Value = 'Cabinet Installation'
$ws.Range("J54").Value = 'Cage Installation'
$ws.Range("J55").Value = 'Cross Connect Single-Mode Fiber'
$ws.Range("J56").Value = 'Equinix Precision Time Standard NTP'
